$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the float formatting pattern "%.2f" to the Pattern row for the
# num_CAL1 (M) and num_2 (P) FLOAT columns.
$ws.Range("M6").Value = "%.2f"
$ws.Range("P6").Value = "%.2f"

# Update the active selection to match the author's final cursor position.
$ws.Range("P6").Select()
